$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells to update: Price (D), Volume(1h) (E) and Hora (G) columns.
# Source cells are stored as literal text (t="inlineStr"), so each
# target cell is pre-formatted as Text (NumberFormat "@") before its
# value is assigned -- otherwise Excel auto-converts numeric- and
# percent-looking strings (e.g. "273.63", "-1.92%") into numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '273.63'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-1.92%'
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = '4'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '26.53'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-3.07%'
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = '4'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '4.873'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '1.24%'
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = '4'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06314'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '0.67%'
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = '4'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.899'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '0.57%'
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = '4'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.351'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '2.74%'
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = '4'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.247'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '31.55%'
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '4'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8723'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-0.60%'
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '4'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1454'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '0.08%'
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '4'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.05147'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '0.42%'
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '4'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07327'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '0.61%'
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '4'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03041'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-3.52%'
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '4'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09034'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.08%'
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '4'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001574'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '1.73%'
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '4'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0006317'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.85%'
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '4'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.006022'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '1.44%'
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '4'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-0.20%'
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '4'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.283'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.33%'
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '4'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '2.55%'
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '4'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1324'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '1.04%'
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '4'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.898'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '1.20%'
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '4'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04419'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '2.15%'
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = '4'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001179'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.51%'
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '4'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004407'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '2.97%'
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = '4'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001199'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '0.05%'
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '4'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001701'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '0.90%'
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '4'
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '4'
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = '4'
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '4'
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '4'
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = '4'
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = '4'
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = '4'
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '4'
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '4'
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = '4'
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = '4'
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = '4'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04026'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-0.41%'
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = '4'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006685'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-2.34%'
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = '4'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1164'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '1.06%'
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = '4'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002109'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-4.47%'
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = '4'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01254'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-11.18%'
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = '4'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005320'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '3.13%'
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = '4'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '3.93%'
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = '4'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.02000'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-33.06%'
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = '4'
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = '4'
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = '4'
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = '4'
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = '4'
